# LCLS/power_calculation.xlsx
# Commit: "Update: m2_q optimization for 9481eV from 100% to 50% power"
#
# 1) Rename the "bw factor" tab to "9482-2DCM".
# 2) Insert a brand-new worksheet "9481-Zig" right before "17795" and fill it
#    with the recalculated 111-220 (9481 eV / 50% power) bandwidth-factor
#    table.
# 3) Restore the cell-cursor / selection state that the author left on the
#    sheets that were actually visited in that session.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Rename "bw factor" -> "9482-2DCM"
# ------------------------------------------------------------------
$bw = $wb.Worksheets.Item("bw factor")
$bw.Name = "9482-2DCM"

# ------------------------------------------------------------------
# 2) Insert the new "9481-Zig" sheet right before "17795"
# ------------------------------------------------------------------
$seventeen = $wb.Worksheets.Item("17795")
$zig = $wb.Worksheets.Add($seventeen)
$zig.Name = "9481-Zig"

# Header row
$zig.Range("A1").Value = "111-220"
$zig.Range("B1").Value = "bragg"
$zig.Range("C1").Value = "asymmetry"
$zig.Range("D1").Value = "b-factor"
$zig.Range("E1").Value = "bandwidth in"
$zig.Range("F1").Value = "bandwidth out"
$zig.Range("G1").Value = "power in"
$zig.Range("H1").Value = "power out"
$zig.Range("I1").Value = "power absorb"

# Row 2 - HHLM1
$zig.Range("A2").Value = "HHLM1"
$zig.Range("B2").Value = 12.042523880204399
$zig.Range("C2").Value = 9
$zig.Range("D2").Formula = "=SIN(RADIANS(B2-C2))/SIN(RADIANS(B2+C2))"
$zig.Range("E2").Formula = "=9481/1000"
$zig.Range("F2").Value = 3.35
$zig.Range("G2").Formula = "=50*0.95^2"
$zig.Range("H2").Formula = "=G2*F2/E2*0.95"
$zig.Range("I2").Formula = "=G2-H2"

# Row 3 - HHLM2
$zig.Range("A3").Value = "HHLM2"
$zig.Range("B3").Value = 19.9159058696595
$zig.Range("C3").Value = 16.9
$zig.Range("D3").Formula = "=SIN(RADIANS(B3-C3))/SIN(RADIANS(B3+C3))"
$zig.Range("E3").Formula = "=F2"
$zig.Range("F3").Value = 3.355
$zig.Range("G3").Formula = "=H2"
$zig.Range("H3").Formula = "=G3*F3/E3*0.95"
$zig.Range("I3").Formula = "=G3-H3"

# Row 4 - HHLM3
$zig.Range("A4").Value = "HHLM3"
$zig.Range("B4").Formula = "=B3"
$zig.Range("C4").Formula = "=-C3"
$zig.Range("D4").Formula = "=SIN(RADIANS(B4-C4))/SIN(RADIANS(B4+C4))"
$zig.Range("E4").Formula = "=F3"
$zig.Range("F4").Value = 3.35
$zig.Range("G4").Formula = "=H3"
$zig.Range("H4").Formula = "=G4*F4/E4*0.95"
$zig.Range("I4").Formula = "=G4-H4"

# Row 5 - HHLM4
$zig.Range("A5").Value = "HHLM4"
$zig.Range("B5").Formula = "=B2"
$zig.Range("C5").Formula = "=-C2"
$zig.Range("D5").Formula = "=SIN(RADIANS(B5-C5))/SIN(RADIANS(B5+C5))"
$zig.Range("E5").Formula = "=F3"
$zig.Range("F5").Value = 3.092
$zig.Range("G5").Formula = "=H3"
$zig.Range("H5").Formula = "=G5*F5/E5*0.95"
$zig.Range("I5").Formula = "=G5-H5"

# Row 6 - C1
$zig.Range("A6").Value = "C1"
$zig.Range("B6").Value = 42.928413263468897
$zig.Range("C6").Value = 0
$zig.Range("D6").Formula = "=SIN(RADIANS(B6-C6))/SIN(RADIANS(B6+C6))"
$zig.Range("E6").Formula = "=F3"
$zig.Range("F6").Value = 0.293
$zig.Range("G6").Formula = "=H3"
$zig.Range("H6").Formula = "=G6*F6/E6*0.95"
$zig.Range("I6").Formula = "=G6-H6"

# Row 7 - C2
$zig.Range("A7").Value = "C2"
$zig.Range("B7").Formula = "=B6"
$zig.Range("C7").Value = -15
$zig.Range("D7").Formula = "=SIN(RADIANS(B7-C7))/SIN(RADIANS(B7+C7))"
$zig.Range("E7").Formula = "=F3"
$zig.Range("F7").Value = 0.174
$zig.Range("G7").Formula = "=H3"
$zig.Range("H7").Formula = "=G7*F7/E7*0.95"
$zig.Range("I7").Formula = "=G7-H7"

# ------------------------------------------------------------------
# 3) Restore the per-sheet selections left by the author
# ------------------------------------------------------------------

# SASE-HHLM-66meV: was topLeftCell A9 / C34 -> now D21
$s = $wb.Worksheets.Item("SASE-HHLM-66meV")
$s.Activate()
$excel.ActiveWindow.Zoom = 100
$s.Range("D21").Select()

# HXRSS-HHLM-66meV: was topLeftCell A4 / H9 -> stays H9 (just drops topLeftCell)
$s = $wb.Worksheets.Item("HXRSS-HHLM-66meV")
$s.Activate()
$excel.ActiveWindow.Zoom = 100
$s.Range("H9").Select()

# Option-2: D19 -> A76
$s = $wb.Worksheets.Item("Option-2")
$s.Activate()
$excel.ActiveWindow.Zoom = 100
$s.Range("A76").Select()

# zig-zag_power_calc: E74 -> D12
$s = $wb.Worksheets.Item("zig-zag_power_calc")
$s.Range("D12").Select()

# 9482-2DCM (was "bw factor"): A14:I20 -> A64:I70
$s = $wb.Worksheets.Item("9482-2DCM")
$s.Range("A64:I70").Select()

# 9481-Zig ends up the active tab, cursor on J8
$zig.Range("J8").Select()
